$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257, shifting the existing data (rows 257-379)
# down to rows 258-380.
$ws.Rows("257:257").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A257").Value = 11
$ws.Range("B257").Value = "Vega Monumental Concepción"
$ws.Range("C257").Value = "Bíobío"
$ws.Range("D257").Value = 44523
$ws.Range("E257").Value = 8
$ws.Range("F257").Value = "Fruta"
$ws.Range("G257").Value = 100108
$ws.Range("H257").Value = "Tropicales y subtropicales"
$ws.Range("I257").Value = 100108006
$ws.Range("J257").Value = "Plátano"
$ws.Range("K257").Value = "Sin especificar"
$ws.Range("L257").Value = "Pintón"
$ws.Range("M257").Value = 650
$ws.Range("N257").Value = 20000
$ws.Range("O257").Value = 21000
$ws.Range("P257").Value = 20538
$ws.Range("Q257").Value = "$/caja 20 kilos"
$ws.Range("R257").Value = "Ecuador"
$ws.Range("S257").Value = 1027
$ws.Range("T257").Value = 20
